$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.329.02'
$ws.Range("E2").Value = '  +10.26%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.576.29'
$ws.Range("E3").Value = '  +11.15%  '

$ws.Range("E4").Value = '  +0.30%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '511.09'
$ws.Range("E5").Value = '  +8.53%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '158.75'
$ws.Range("E6").Value = '  +11.07%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.991'
$ws.Range("E7").Value = '  -1.05%  '

$ws.Range("E8").Value = '  +0.93%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.633.41'
$ws.Range("E9").Value = '  +13.86%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.14'
$ws.Range("E10").Value = '  +13.77%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.105'
$ws.Range("E11").Value = '  +10.94%  '

$ws.Range("E12").Value = '  +7.55%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.031.10'
$ws.Range("E14").Value = '  +11.02%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '60.063.97'
$ws.Range("E15").Value = '  +9.49%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '22.14'
$ws.Range("E16").Value = '  +12.09%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000141'
$ws.Range("E17").Value = '  +10.21%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.618.02'
$ws.Range("E18").Value = '  +12.29%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.82'
$ws.Range("E19").Value = '  +7.31%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '343.82'
$ws.Range("E20").Value = '  +10.68%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.56'
$ws.Range("E21").Value = '  +11.81%  '

$ws.Range("E22").Value = '  +10.52%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("E23").Value = '  -0.15%  '

$ws.Range("E24").Value = '  +8.22%  '

$ws.Range("E25").Value = '  +8.85%  '

$ws.Range("E26").Value = '  +14.14%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.694.82'
$ws.Range("E27").Value = '  +10.61%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.989'
$ws.Range("E28").Value = '  -1.29%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0₃0863'
$ws.Range("E29").Value = '  +18.45%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.49'
$ws.Range("E30").Value = '  +7.82%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.996'
$ws.Range("E31").Value = '  -0.49%  '

$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '19.67'
$ws.Range("E32").Value = '  +9.86%  '

$ws.Range("B33").Value = 'Monero'
$ws.Range("C33").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '157.30'
$ws.Range("E33").Value = '  +8.46%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.59'
$ws.Range("E34").Value = '  +9.00%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.60'
$ws.Range("E35").Value = '  +11.62%  '

$ws.Range("E36").Value = '  +12.53%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.98'
$ws.Range("E37").Value = '  +12.21%  '

$ws.Range("B38").Value = 'Fetch.AI'
$ws.Range("C38").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.876'
$ws.Range("E38").Value = '  +9.96%  '

$ws.Range("B39").Value = 'Bittensor'
$ws.Range("C39").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '309.86'
$ws.Range("E39").Value = '  +25.79%  '

$ws.Range("E40").Value = '  +13.58%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.77'
$ws.Range("E41").Value = '  +12.91%  '

$ws.Range("E42").Value = '  +6.37%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.638'
$ws.Range("E43").Value = '  +12.41%  '

$ws.Range("B44").Value = 'SuiNetwork'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.800'
$ws.Range("E44").Value = '  +31.22%  '

$ws.Range("B45").Value = 'Hedera'
$ws.Range("C45").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0578'
$ws.Range("E45").Value = '  +13.49%  '

$ws.Range("B46").Value = 'Stellar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.102'
$ws.Range("E46").Value = '  +0.95%  '

$ws.Range("B47").Value = 'RenderToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '5.05'
$ws.Range("E47").Value = '  +17.34%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '19.92'
$ws.Range("E48").Value = '  +21.64%  '

$ws.Range("B49").Value = 'FirstDigitalUSD'
$ws.Range("C49").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.989'
$ws.Range("E49").Value = '  -0.94%  '

$ws.Range("E50").Value = '  +9.05%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.004.98'
$ws.Range("E51").Value = '  +13.70%  '
